{"js": "// Applies the \"e2, l5, s1 printp=out changes\" edit set to the technical\n// specification document:\n//   - Enquiry date bumped forward\n//   - Project value corrected\n//   - Cooling-circuit \"Heat Rejected\" + direct-fired \"Heat Input\" values updated\n//   - Fuel switched from HSD (liquid, kcal/kg, GCV value, 9.3) to Natural Gas\n//     (kcal/Nm\u00b3, GCV value, Nm\u00b3/hr units)\n//   - a previously-blank \"Gas Pressure\" spec row is populated\n\nconst body = context.document.body;\n\n// ---- helper: replace a single, unique, exact text match in the body ----\nasync function replaceUniqueText(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---- helper: replace the Nth (0-based) occurrence of a text match ----\nasync function replaceNthOccurrence(searchText, occurrenceIndex, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length <= occurrenceIndex) {\n    throw new Error(\n      `Expected at least ${occurrenceIndex + 1} matches for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n  results.items[occurrenceIndex].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---- helper: set the text of a (currently empty) table cell, preserving\n//      the existing paragraph/run formatting instead of minting a new run ----\nasync function setCellText(table, rowIndex, colIndex, text) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.getRange().insertText(text, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Enquiry date/time value\nawait replaceUniqueText(\"10-Sep-2021, 12:12\", \"20-Sep-2021, 13:11\");\n\n// 2) Project value \"aa\" -> \"a\"\nawait replaceUniqueText(\"aa\", \"a\");\n\n// 3) Cooling Water Circuit - Heat Rejected value\nawait replaceUniqueText(\"213349.4\", \"213366.8\");\n\n// 4) Direct Fired Circuit - Heat Input value\nawait replaceUniqueText(\"98060.7\", \"98076.6\");\n\n// 5) Fuel Type\nawait replaceUniqueText(\"HSD\", \"NaturalGas\");\n\n// 6) Calorific Value unit\nawait replaceUniqueText(\"kcal/kg\", \"kcal/Nm\\u00B3\");\n\n// 7) Calorific Value value\nawait replaceUniqueText(\"10200\", \"9000\");\n\n// 8) Fuel consumption unit - \"GCV\" appears twice (Calorific value type row,\n//    then Fuel consumption row); the second occurrence is the one that\n//    switches to a volumetric flow unit.\nawait replaceNthOccurrence(\"GCV\", 1, \"Nm\\u00B3/hr\");\n\n// 9) Fuel consumption value\nawait replaceUniqueText(\"9.3\", \"10.5\");\n\n// 10-13) Populate the previously blank \"Gas Pressure\" spec row. Locate it\n// dynamically via the preceding row's label so the edit is resilient to\n// any row-index drift.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst directFiredTable = tables.items[1];\ndirectFiredTable.load(\"values\");\nawait context.sync();\n\nlet gasPressureRow = -1;\nfor (let i = 0; i < directFiredTable.values.length; i++) {\n  if (directFiredTable.values[i][1] === \"Exhaust Gas duct size\") {\n    gasPressureRow = i + 1;\n    break;\n  }\n}\nif (gasPressureRow === -1) {\n  throw new Error('Could not locate the row following \"Exhaust Gas duct size\"');\n}\n\nawait setCellText(directFiredTable, gasPressureRow, 0, \"7.\");\nawait setCellText(directFiredTable, gasPressureRow, 1, \"Gas Pressure\");\nawait setCellText(directFiredTable, gasPressureRow, 2, \"mbar\");\nawait setCellText(directFiredTable, gasPressureRow, 3, \"100\");\n", "ps1": "# Applies the \"e2, l5, s1 printp=out changes\" edit set to the technical\n# specification document:\n#   - Enquiry date bumped forward\n#   - Project value corrected\n#   - Cooling-circuit \"Heat Rejected\" + direct-fired \"Heat Input\" values updated\n#   - Fuel switched from HSD (liquid, kcal/kg, GCV value, 9.3) to Natural Gas\n#     (kcal/Nm3, GCV value, Nm3/hr units)\n#   - a previously-blank \"Gas Pressure\" spec row is populated\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for '$findText'\"\n    }\n}\n\n# 1) Enquiry date/time value\nReplace-UniqueText \"10-Sep-2021, 12:12\" \"20-Sep-2021, 13:11\"\n\n# 2) Project value \"aa\" -> \"a\"\nReplace-UniqueText \"aa\" \"a\"\n\n# 3) Cooling Water Circuit - Heat Rejected value\nReplace-UniqueText \"213349.4\" \"213366.8\"\n\n# 4) Direct Fired Circuit - Heat Input value\nReplace-UniqueText \"98060.7\" \"98076.6\"\n\n# 5) Fuel Type\nReplace-UniqueText \"HSD\" \"NaturalGas\"\n\n# 6) Calorific Value unit\nReplace-UniqueText \"kcal/kg\" \"kcal/Nm\u00b3\"\n\n# 7) Calorific Value value\nReplace-UniqueText \"10200\" \"9000\"\n\n# 8) Fuel consumption unit - \"GCV\" appears twice in the document (once for\n# the \"Calorific value type\" row, once for the \"Fuel consumption\" row).\n# wdReplaceAll would touch both, so scope a range to just after \"Fuel\n# consumption\" and replace only the one match inside it.\n$anchor = $d.Content\n$anchor.Start = 0\n$anchorFind = $anchor.Find\n$anchorFind.ClearFormatting()\n$anchorFind.Text = \"Fuel consumption\"\n$anchorFound = $anchorFind.Execute($anchorFind.Text, $false, $false, $false, $false, $false, $true, 0, $false, \"\", 0)\nif (-not $anchorFound) {\n    throw \"Could not locate 'Fuel consumption' anchor\"\n}\n$afterRange = $d.Range($anchor.End, $d.Content.End)\n$gcvFind = $afterRange.Find\n$gcvFind.ClearFormatting()\n$gcvFind.Replacement.ClearFormatting()\n$gcvFind.Text = \"GCV\"\n$gcvFind.Replacement.Text = \"Nm\u00b3/hr\"\n$gcvResult = $gcvFind.Execute($gcvFind.Text, $false, $false, $false, $false, $false, $true, 0, $false, $gcvFind.Replacement.Text, 1)\nif (-not $gcvResult) {\n    throw \"Could not replace the Fuel consumption unit (GCV)\"\n}\n\n# 9) Fuel consumption value\nReplace-UniqueText \"9.3\" \"10.5\"\n\n# 10-13) Populate the previously blank \"Gas Pressure\" spec row. Locate it\n# dynamically via the preceding row's label so the edit is resilient to any\n# row-index drift.\n$table = $d.Tables.Item(2)\n$gasPressureRow = -1\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $labelText = $table.Cell($r, 2).Range.Text\n    if ($labelText -like \"*Exhaust Gas duct size*\") {\n        $gasPressureRow = $r + 1\n        break\n    }\n}\nif ($gasPressureRow -eq -1) {\n    throw \"Could not locate the row following 'Exhaust Gas duct size'\"\n}\n\n$table.Cell($gasPressureRow, 1).Range.Text = \"7.\"\n$table.Cell($gasPressureRow, 2).Range.Text = \"Gas Pressure\"\n$table.Cell($gasPressureRow, 3).Range.Text = \"mbar\"\n$table.Cell($gasPressureRow, 4).Range.Text = \"100\"\n"}
